$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-20 Thursday" "2025-03-21 Friday"
Replace-Text "543÷8=67, 7" "931÷5=186, 1"
Replace-Text "862÷6=143, 4" "384÷9=42, 6"
Replace-Text "437÷5=87, 2" "119÷5=23, 4"
Replace-Text "767÷6=127, 5" "880÷9=97, 7"
Replace-Text "329÷8=41, 1" "423÷6=70, 3"
Replace-Text "511÷2=255, 1" "553÷8=69, 1"
Replace-Text "753÷9=83, 6" "284÷3=94, 2"
Replace-Text "869÷6=144, 5" "939÷9=104, 3"
Replace-Text "204÷4=51, 0" "156÷6=26, 0"
Replace-Text "280÷3=93, 1" "352÷7=50, 2"
Replace-Text "813÷2=406, 1" "833÷3=277, 2"
Replace-Text "911÷2=455, 1" "970÷9=107, 7"
Replace-Text "707÷8=88, 3" "636÷7=90, 6"
Replace-Text "351÷6=58, 3" "449÷9=49, 8"
Replace-Text "535÷5=107, 0" "746÷7=106, 4"
Replace-Text "965÷6=160, 5" "187÷5=37, 2"
Replace-Text "527÷5=105, 2" "644÷6=107, 2"
Replace-Text "747÷4=186, 3" "682÷9=75, 7"
Replace-Text "934÷3=311, 1" "337÷7=48, 1"
Replace-Text "875÷4=218, 3" "905÷8=113, 1"
Replace-Text "554÷4=138, 2" "292÷8=36, 4"
Replace-Text "509÷2=254, 1" "831÷6=138, 3"
Replace-Text "686÷4=171, 2" "989÷9=109, 8"
Replace-Text "532÷8=66, 4" "694÷6=115, 4"
Replace-Text "407÷5=81, 2" "716÷9=79, 5"
